$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Move the "correct" label from D2 to E2, and add a new header in E1
$ws.Range("D2").Clear()
$ws.Range("E2").Value = "correct"
$ws.Range("E1").Value = "reference MIC"

# Highlight the whole "reference MIC" column (E) with a light blue fill
# (Blue, Accent 1, Lighter 80%)
$ws.Columns.Item(5).Interior.Color = 15983578

# Move the active selection to G13, matching the latest view state
$ws.Range("G13").Select()
